$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new price-report date (2022-02-23, serial 44615) was
# added for "Terminal La Palmera de La Serena - Palta". It lands at the
# top of this product's block (row 647), so insert 3 blank rows there
# and push the existing data (old rows 647:671) down to 650:674.
$ws.Rows("647:649").Insert()

# Especial / Primera / Segunda rows for the new date, same shape as the
# other "caja de 17 kilos" / "Provincia de Limarí" entries in this block.
$newRows = @(
    @("Especial", 360, 2500, 2600, 2550, 2550),
    @("Primera",  240, 2200, 2300, 2250, 2250),
    @("Segunda",  200, 1900, 2000, 1950, 1950)
)

for ($i = 0; $i -lt 3; $i++) {
    $r = 647 + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44615
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100106
    $ws.Cells.Item($r, 8).Value = "Oleaginosos"
    $ws.Cells.Item($r, 9).Value = 100106002
    $ws.Cells.Item($r, 10).Value = "Palta"
    $ws.Cells.Item($r, 11).Value = "Hass"
    $ws.Cells.Item($r, 12).Value = $vals[0]
    $ws.Cells.Item($r, 13).Value = $vals[1]
    $ws.Cells.Item($r, 14).Value = $vals[2]
    $ws.Cells.Item($r, 15).Value = $vals[3]
    $ws.Cells.Item($r, 16).Value = $vals[4]
    $ws.Cells.Item($r, 17).Value = "$/kilo (en caja de 17 kilos)"
    $ws.Cells.Item($r, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 19).Value = $vals[5]
    $ws.Cells.Item($r, 20).Value = 1
}
